# Update average_county_temperature (column AD) with NOAA data for
# facilities under NAICS 311615. Each facility/site shares the same
# updated temperature value across its contiguous block of unit rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2:AD10").Value = 19.30324074074072
$ws.Range("AD11:AD16").Value = 17.25771604938272
$ws.Range("AD17:AD25").Value = 13.62268518518517
$ws.Range("AD65:AD67").Value = 13.75752314814816
$ws.Range("AD68:AD73").Value = 19.79629629629628
$ws.Range("AD74:AD79").Value = 0.8611111111111096
$ws.Range("AD83:AD88").Value = 0.8611111111111096
$ws.Range("AD92:AD94").Value = 5.486111111111112
$ws.Range("AD95:AD100").Value = 16.86342592592595
$ws.Range("AD101:AD106").Value = 5.486111111111112
$ws.Range("AD107:AD109").Value = 12.41429539295394
$ws.Range("AD113:AD115").Value = 19.36574074074073
$ws.Range("AD140:AD142").Value = 12.41429539295394
$ws.Range("AD146:AD154").Value = 12.41429539295394
